# Updates the NBA "team leaders" workbook: refreshes each leaderboard sheet
# (rank / team name / value) with the latest stats.
$wb = $excel.ActiveWorkbook

$data = @{
    "Arremessos %" = @(
        @(1, "Oklahoma City Thunder", "50.9"),
        @(2, "Los Angeles Lakers",    "49.0"),
        @(3, "Miami Heat",            "48.1"),
        @(4, "Minnesota Timberwolves","47.6"),
        @(5, "Indiana Pacers",        "47.5")
    )
    "Diferencial de Pontos" = @(
        @(1, "Minnesota Timberwolves","+18.5"),
        @(2, "Oklahoma City Thunder", "+17.0"),
        @(3, "Denver Nuggets",        "+6.6"),
        @(4, "Boston Celtics",        "+5.0"),
        @(5, "Orlando Magic",         "+4.6")
    )
    "Pontos" = @(
        @(1, "Minnesota Timberwolves","112.5"),
        @(2, "Philadelphia 76ers",    "110.0"),
        @(3, "New York Knicks",       "109.7"),
        @(4, "Indiana Pacers",        "109.5"),
        @(5, "Denver Nuggets",        "109.0")
    )
    "Pontos Permitidos" = @(
        @(1, "Orlando Magic",          "92.0"),
        @(1, "Oklahoma City Thunder",  "92.0"),
        @(3, "Minnesota Timberwolves", "94.0"),
        @(4, "LA Clippers",            "96.5"),
        @(5, "Cleveland Cavaliers",    "96.7")
    )
    "Rebotes" = @(
        @(1, "Denver Nuggets",         "48.3"),
        @(2, "Indiana Pacers",         "47.5"),
        @(3, "LA Clippers",            "46.0"),
        @(4, "Dallas Mavericks",       "45.5"),
        @(4, "Minnesota Timberwolves", "45.5")
    )
    "Tocos" = @(
        @(1, "Dallas Mavericks",      "8.0"),
        @(2, "Philadelphia 76ers",    "7.0"),
        @(2, "Boston Celtics",        "7.0"),
        @(4, "Cleveland Cavaliers",   "6.0"),
        @(5, "Indiana Pacers",        "5.5")
    )
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $data[$sheetName]
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $i + 2
        $row = $rows[$i]
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        # Leading apostrophe forces the numeric-looking stat to stay text,
        # matching the original "Valor" column (e.g. "50.9", "+18.5").
        $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    }
}
